$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.245.94"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "2.364.20"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.685"
$ws.Range("E5").Value = "  +4.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.26"
$ws.Range("E6").Value = "  +3.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.43"
$ws.Range("E7").Value = "  +3.49%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  +27.65%  "
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.18"
$ws.Range("E11").Value = "  +20.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.51"
$ws.Range("E12").Value = "  +19.93%  "
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").Value = "2.713.35"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.96"
$ws.Range("E15").Value = "  +6.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.917"
$ws.Range("E16").Value = "  +7.57%  "
$ws.Range("D17").Value = "2.357.56"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "44.241.99"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("E20").Value = "  +5.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "78.56"
$ws.Range("E21").Value = "  +5.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "256.48"
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.60"
$ws.Range("E23").Value = "  +4.55%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.75"
$ws.Range("E25").Value = "  -4.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.79"
$ws.Range("E27").Value = "  +5.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.66"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.63"
$ws.Range("E29").Value = "  +5.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.32"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("E32").Value = "  +4.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.43"
$ws.Range("E33").Value = "  +8.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0757"
$ws.Range("E34").Value = "  +9.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.39"
$ws.Range("E35").Value = "  +5.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.91"
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0274"
$ws.Range("E39").Value = "  +6.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.33"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.12"
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.197"
$ws.Range("E43").Value = "  +15.94%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  +12.40%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.26"
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("E46").Value = "  +4.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.08"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.48"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").Value = "1.458.27"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("E51").Value = "  +2.93%  "
